# Append two new slides ("Summary" and "Contact Information") to the
# end of the deck, cloned from the last existing slide (slide 18,
# "Predicted Salary by Age") so they inherit its decorative shapes and
# title placeholder, then strip the picture and retitle each one.

$p = $ppt.ActivePresentation

# --- New slide 1: "Summary" (becomes sldId 304, slide index 19) -----
$lastSlide = $p.Slides.Item($p.Slides.Count)
$summarySlide = $lastSlide.Duplicate()
$summarySlide = $p.Slides.Item($summarySlide.SlideIndex)

$summarySlide.Shapes.Item("Picture 9").Delete()
$summarySlide.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Summary"

# --- New slide 2: "Contact Information" (becomes sldId 305, slide 20)
$contactSlide = $summarySlide.Duplicate()
$contactSlide = $p.Slides.Item($contactSlide.SlideIndex)

$contactSlide.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Contact Information"
